# Applies the edits described by the commit:
#   - Fill in the "등급" (grade) column (J5:J12) on sheet "제 1작업" with an
#     IF/AND formula classifying each row as "A" or "B".
#   - Fill in the two summary formulas in row 13/14 that were left blank
#     before (E13 한식 count, J13 min menu count, E14 DSUM total, J14
#     VLOOKUP lookup).
#   - Widen column E on that sheet.
#   - Add conditional formatting (bold, blue font) to B5:J12 for rows whose
#     전월배달건수 (H) is below 300.
#   - Restore the cell-cursor position on three sheets to match where the
#     author left it.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("제 1작업")

# --- New "등급" formulas in column J (rows 5-12) -----------------------
$ws1.Range("J5").Formula  = '=IF(AND(F5>=15,H5>=300),"A","B")'
$ws1.Range("J6:J12").Formula = '=IF(AND(F6>=15,H6>=300),"A","B")'

# --- Summary formulas that used to be blank -----------------------------
$ws1.Range("E13").Formula = '=COUNTIF(D5:D12,"한식")&"개"'
$ws1.Range("J13").Formula = '=MIN(메뉴수)'
$ws1.Range("E14").Formula = '=DSUM(B4:H12,H4,D4:D5)'
$ws1.Range("J14").Formula = '=VLOOKUP(H14,B5:H12,7,FALSE)'

# --- Column E is now a little wider -------------------------------------
$ws1.Columns("E").ColumnWidth = 13.25

# --- Conditional formatting: bold blue font when 전월배달건수 < 300 -----
$cfRange = $ws1.Range("B5:J12")
$cf = $cfRange.FormatConditions.Add(2, 0, '=$H5<300')
$cf.Font.Bold = $true
$cf.Font.Italic = $false
$cf.Font.Color = 12611584

# --- Leave the selection where the author left it on each sheet --------
$ws1.Activate()
$ws1.Range("J6").Select()

$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Activate()
$ws5.Range("E12").Select()

$ws6 = $wb.Worksheets.Item("Sheet1")
$ws6.Activate()
$ws6.Range("E13").Select()

$ws1.Activate()
